$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assurance Qualité")
Write-Host $ws.Name
